$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.613.13"
$ws.Range("E2").Value = "  +5.91%  "
$ws.Range("D3").Value = "3.489.37"
$ws.Range("E3").Value = "  +6.37%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.77"
$ws.Range("E5").Value = "  +8.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "554.02"
$ws.Range("E6").Value = "  +6.69%  "
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").Value = "3.485.29"
$ws.Range("E8").Value = "  +6.27%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.638"
$ws.Range("E10").Value = "  +6.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.78"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.149"
$ws.Range("E12").Value = "  +12.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +7.41%  "
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "4.048.66"
$ws.Range("E15").Value = "  +6.58%  "
$ws.Range("D16").Value = "3.494.43"
$ws.Range("E16").Value = "  +6.64%  "
$ws.Range("D17").Value = "68.182.19"
$ws.Range("E17").Value = "  +6.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("E18").Value = "  +4.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.29"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  +7.77%  "
$ws.Range("E21").Value = "  +6.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.95"
$ws.Range("E22").Value = "  +9.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.13"
$ws.Range("E23").Value = "  +10.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.96"
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.74"
$ws.Range("E25").Value = "  +6.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.21"
$ws.Range("E26").Value = "  +8.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("E27").Value = "  +11.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.25"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.85"
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.64"
$ws.Range("E30").Value = "  +5.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.35"
$ws.Range("E31").Value = "  +6.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "684.68"
$ws.Range("E32").Value = "  +7.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.92"
$ws.Range("E33").Value = "  +5.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.70"
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("E35").Value = "  +6.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.61"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "0.0₃0834"
$ws.Range("E37").Value = "  +21.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.93"
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  +27.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  +16.45%  "
$ws.Range("E43").Value = "  +10.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "3.058.97"
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("E46").Value = "  +12.16%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0422"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.31"
$ws.Range("E48").Value = "  +12.74%  "
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.01"
$ws.Range("E50").Value = "  +16.53%  "
$ws.Range("E51").Value = "  +4.89%  "
